# Weekly update: a new "Ajo" price record (week of 2022-01-05) is added to
# the Macroferia Regional de Talca dataset. The new record is inserted as
# row 123 (in Fecha/date order used by the source feed), which pushes the
# existing rows 123-225 down to 124-226, growing the sheet by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 123, shifting rows 123:225 down to 124:226.
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row with the new weekly record. Columns A, B,
# C, E, F, G, H, I, O, Q, R keep the same series values (same market,
# region, category, variety, quality, origin, classification) as the rest
# of this "Ajo" subset; D (Fecha), J (Volumen), K/L/M (precios), N (unidad
# de comercializacion) and P (Precio $/Kg) hold the new observation.
$ws.Cells.Item(123, 1).Value = 5
$ws.Cells.Item(123, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(123, 3).Value = "Maule"
$ws.Cells.Item(123, 4).Value = 44566
$ws.Cells.Item(123, 5).Value = 7
$ws.Cells.Item(123, 6).Value = 100112003
$ws.Cells.Item(123, 7).Value = "Ajo"
$ws.Cells.Item(123, 8).Value = "Chino"
$ws.Cells.Item(123, 9).Value = "Primera"
$ws.Cells.Item(123, 10).Value = 200
$ws.Cells.Item(123, 11).Value = 20000
$ws.Cells.Item(123, 12).Value = 20000
$ws.Cells.Item(123, 13).Value = 20000
$ws.Cells.Item(123, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(123, 15).Value = "China"
$ws.Cells.Item(123, 16).Value = 2000
$ws.Cells.Item(123, 17).Value = 10
$ws.Cells.Item(123, 18).Value = "Hortaliza"

Write-Output "Inserted new row 123; sheet now spans to row $($ws.UsedRange.Rows.Count)."
